$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Update the header text to upper-case the month name.
    $ws.Range("A1").Value = "Month/Year: DECEMBER 2021"

    # Clear the tab color (best effort - remove the blue sheet-tab highlight).
    $ws.Tab.ColorIndex = -4142
}
